$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the preprocessing/confounds folder version string (fmriprep 20.0.1 -> 20.0.6)
# These values live in the dataFolder (G) and confoundsFolder (H) columns for rows 2-5
$ws.Range("G2:H5").Value = "preproc_fmriprep-20.0.6"

# Update the active selection to H2 (was H13)
$ws.Range("H2").Select() | Out-Null
